$d = $word.ActiveDocument

function Split-Run($startPos, $endPos) {
    # Forces a run boundary at [startPos, endPos) by toggling Bold on and
    # back off. This does not change the visible text/formatting, it only
    # causes the run-coalescing pass to leave this span as its own run.
    if ($endPos -gt $startPos) {
        $r = $d.Range($startPos, $endPos)
        $r.Bold = 1
        $r.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# Change 1: fix typo "копоновку" -> "компоновку"
# ---------------------------------------------------------------------------

$rng1 = $d.Content
$null = $rng1.Find.Execute(", я передала на обработку компоновщику объектный файл.После этого я сделала копоновку, с помощью команды")
$run1Start = $rng1.Start
$rng1.Text = ", я передала на обработку компоновщику объектный файл.После этого я сделала компоновку, с помощью команды"

# The assignment above merges the run with its predecessor
# ("“ld -m elf_i386 hello.o -o hello”") and with every run that follows it
# through to the end of the paragraph. Re-discover each of those original
# runs (now substrings of one big run) and re-split them back out.

$quote1 = $d.Range($run1Start - 60, $run1Start)
$null = $quote1.Find.Execute("“ld -m elf_i386 hello.o -o hello”")
Split-Run $quote1.Start $quote1.End

$fixedEnd = $d.Range($run1Start, $run1Start)
$null = $fixedEnd.Find.Execute(", я передала на обработку компоновщику объектный файл.После этого я сделала компоновку, с помощью команды")
$afterFixed = $fixedEnd.End

$sp2 = $d.Range($afterFixed, $afterFixed + 1)
Split-Run $sp2.Start $sp2.End

$quote2 = $d.Range($afterFixed, $afterFixed + 60)
$null = $quote2.Find.Execute("“ld -m elf_i386 obj.o -o main”")
Split-Run $quote2.Start $quote2.End

$sp3 = $d.Range($quote2.End, $quote2.End + 1)
Split-Run $sp3.Start $sp3.End

$fig1 = $d.Range($quote2.End, $quote2.End + 30)
$null = $fig1.Find.Execute("(Рисунок 6.1)")
Split-Run $fig1.Start $fig1.End

# ---------------------------------------------------------------------------
# Change 2: extend the "(Рисунок 8.2)" sentence and append the follow-up
# sentence about the printed "Drekina Arina" text.
# ---------------------------------------------------------------------------

$rng2 = $d.Content
$null = $rng2.Find.Execute("запустила исполняемый файл. (Рисунок 8.2)")
$origStart = $rng2.Start
$rng2.Text = "запустила исполняемый файл. (Рисунок 8.2). На монитор вывелся измененный текст"

# This merges in the single space-only run that used to sit right before
# this run. Split it back off.
Split-Run ($origStart - 1) $origStart

# Append the three new runs after the (now corrected) run, using
# InsertAfter on a collapsed range so no existing run is disturbed.
$tail = $d.Content
$null = $tail.Find.Execute("запустила исполняемый файл. (Рисунок 8.2). На монитор вывелся измененный текст")
$tail.Collapse(0)
$tail.InsertAfter(" ")
$tail.Collapse(0)
$tail.InsertAfter("“Drekina Arina”")
$tail.Collapse(0)
$tail.InsertAfter(", это означает, что во время редактирования я не допустила ошибок.")
